$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.730.46"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -5.44%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.964.79"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -6.91%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.66"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.52%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "124.61"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -8.01%  "

# Row 7
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.957.05"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -7.07%  "

# Row 9
$ws.Range("E9").Value = "  -2.86%  "

# Row 10
$ws.Range("E10").Value = "  -6.18%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.06"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -3.46%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.437"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.73%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000221"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -6.77%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.38"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -6.84%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.118"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.74%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.455.54"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -6.91%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.964.16"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -7.04%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "59.764.47"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.31%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.16"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -6.34%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "432.28"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -6.60%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.03"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -7.12%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.658"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.68%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.98"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -8.76%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.69"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.26%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "78.93"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.75%  "

# Row 26
$ws.Range("E26").Value = "  -0.07%  "

# Row 27
$ws.Range("E27").Value = "  -0.16%  "

# Row 28
$ws.Range("E28").Value = "  -5.97%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.17"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -7.07%  "

# Row 30
$ws.Range("E30").Value = "  -7.99%  "

# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.28"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -7.30%  "

# Row 32
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.10"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -9.88%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0927"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -9.83%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.16"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -9.29%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.944"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -8.78%  "

# Row 36
$ws.Range("E36").Value = "  -4.33%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "49.52"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.54%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0655"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -7.70%  "

# Row 39
$ws.Range("B39").Value = "Cosmos"
$ws.Range("C39").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.89"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.55%  "

# Row 40
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0356"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -8.58%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.108"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -3.64%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "378.84"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -6.90%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.46"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -7.17%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.618.78"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -6.78%  "

# Row 45
$ws.Range("E45").Value = "  +0.06%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.235"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -7.16%  "

# Row 47
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.24"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.09%  "

# Row 48
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.97"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -7.39%  "

# Row 49
$ws.Range("E49").Value = "  -4.80%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.26"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -8.07%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "31.07"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -11.48%  "
